$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Week 3 notes for Team 4 (row 6)
$ws.Range("N6").Value = "Team has investigated various datasets and methods currently used in research. The team has permission to use these datasets for research. So, they are currently working on exploratory data analysis on these datasets to decide which one to use."

# Week 3 notes for Team 5 (row 7)
$ws.Range("N7").Value = "Data analysis has been completed. The team is working on the KNN model. They are also checking to see if they can use any other clustering technique or augment their dataset using summaries of movie from a different data source."

# Week 3 notes for Team 6 (row 8)
$ws.Range("N8").Value = "Data pre-processing has been started. Random Forest and Stacked Classifier model is being planned."

# Week 3 notes for Team 7 (row 9)
$ws.Range("N9").Value = "Data pre-processing has been completed. The team is currently in the exploratory data analysis and feature engineering phase."

# Row 8 previously held only a one-line comment; the new note wraps onto
# multiple lines so the row needs to grow to show it fully.
$ws.Rows.Item(8).RowHeight = 37.3

# Move the view / selection to where the editor ended up after typing the
# last note.
$ws.Application.Goto($ws.Range("J2"))
$ws.Range("N2").Select()
